# Water/fertilizer shuffle log - add three new entries (rows 38-40) and
# a new note, then leave the sheet selection where the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- new log rows -----------------------------------------------------
# Row 38: 2025-09-27, watered / soaked, no extra note
$ws.Cells.Item(38, 1).Value = 45927
$ws.Cells.Item(38, 3).Value = "watered"
$ws.Cells.Item(38, 4).Value = "soaked"

# Row 39: 2025-09-28, watered / soaked, no extra note
$ws.Cells.Item(39, 1).Value = 45928
$ws.Cells.Item(39, 3).Value = "watered"
$ws.Cells.Item(39, 4).Value = "soaked"

# Row 40: 2025-10-01, watered / soaked, thanks note
$ws.Cells.Item(40, 1).Value = 45931
$ws.Cells.Item(40, 3).Value = "watered"
$ws.Cells.Item(40, 4).Value = "soaked"
$ws.Cells.Item(40, 5).Value = "thanks Zarra"

# Column B holds the "day of week" helper formula (=A<row>); fill it as
# one range assignment so the new rows share a single formula group,
# same as the existing B3:B37 block.
$ws.Range("B38:B40").Formula = "=A38"

# --- restore the view/selection state the author saved with ----------
# (window geometry isn't surfaced by this host the same way a live
# Excel instance would be, but set it anyway in case it is honoured)
try {
  $win = $wb.Windows.Item(1)
  $win.Left = 2580
  $win.Top = 2580
  $win.Width = 17280
  $win.Height = 9960
} catch {
}

$ws.Range("C34").Select() | Out-Null
